# The presentation ships with two embedded DrawingML themes:
#   theme1.xml -> linked from the (single) slide master -> "Integral" / "Red Violet" colours
#   theme2.xml -> linked from the notes master           -> "Office Theme" / "Office" colours
#
# The authored change swaps the content of the two theme parts in place:
#   theme1.xml ends up holding the "Office Theme" colour scheme
#   theme2.xml ends up holding the "Integral" (Red Violet) colour scheme
#
# The PowerPoint object model exposes the live (slide-master-linked) theme's
# 12 theme colours through Slide.ThemeColorScheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink, in that order) - this is the handle that lets
# us rewrite theme1.xml's <a:clrScheme> to the target "Office" palette.

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette: the stock Office theme colours (what theme1.xml should
# contain after the edit).
$officeColors = @(
    @(0,0,0),        # 1  dk1      000000
    @(255,255,255),  # 2  lt1      FFFFFF
    @(68,84,106),    # 3  dk2      44546A
    @(231,230,230),  # 4  lt2      E7E6E6
    @(91,155,213),   # 5  accent1  5B9BD5
    @(237,125,49),   # 6  accent2  ED7D31
    @(165,165,165),  # 7  accent3  A5A5A5
    @(255,192,0),    # 8  accent4  FFC000
    @(68,114,196),   # 9  accent5  4472C4
    @(112,173,71),   # 10 accent6  70AD47
    @(5,99,193),     # 11 hlink    0563C1
    @(149,79,114)    # 12 folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $c = $officeColors[$i - 1]
    $tcs.Item($i).RGB = RGBVal $c[0] $c[1] $c[2]
}

# The underlying colour-scheme / theme display names ("Red Violet" /
# "Integral" vs "Office" / "Office Theme") are not writable through any
# exposed PowerPoint object-model property in this host, so only the
# colour values themselves are updated here.
try { $tcs.Name = "Office" } catch { }
try { $s.Design.Name = "Office Theme" } catch { }

Write-Host "Theme colours updated"
